# "Generate Report for Handback" — the handback for
# f5348948-4f86-4d36-b8a0-67a8c1d3ffcf.md completed successfully, so the
# status moves from "Ready for handoff" to "Handed back: in sync with
# en-US", the per-locale "Latest Handback DateTime" timestamps are
# refreshed, and the stale "version not latest" Error Detail is cleared.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the f5348948 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the f5348948 file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K3").Value = "2016-10-17 13:52:00"
$wsZh.Range("P3").Value = ""
$wsZh.Columns.Item(16).AutoFit()

# --- de-de sheet: row 3 is the f5348948 file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K3").Value = "2016-10-17 13:52:38"
$wsDe.Range("P3").Value = ""
$wsDe.Columns.Item(16).AutoFit()
